# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / LevePriceNQ / LeveProfitNQ
# and their HQ counterparts where applicable) across the per-job profit
# sheets, per the latest market data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 47746.5
$ws.Range("I62").Value = 17249.572
$ws.Range("K62").Value = 17249.572
$ws.Range("M62").Value = -16625.572
$ws.Range("H65").Value = 47746.5
$ws.Range("I65").Value = 17249.572
$ws.Range("K65").Value = 86247.86
$ws.Range("M65").Value = -83127.86
$ws.Range("H80").Value = 20673.16
$ws.Range("I80").Value = 6477.9414
$ws.Range("J80").Value = 50838
$ws.Range("K80").Value = 19433.8242
$ws.Range("L80").Value = 152514
$ws.Range("M80").Value = -18435.8242
$ws.Range("N80").Value = -154510
$ws.Range("H83").Value = 20673.16
$ws.Range("I83").Value = 6477.9414
$ws.Range("J83").Value = 50838
$ws.Range("K83").Value = 58301.47259999999
$ws.Range("L83").Value = 457542
$ws.Range("M83").Value = -53309.47259999999
$ws.Range("N83").Value = -467526
$ws.Range("H111").Value = 20837828
$ws.Range("I111").Value = 41667332
$ws.Range("K111").Value = 125001996
$ws.Range("M111").Value = -124998929
$ws.Range("H116").Value = 14711228
$ws.Range("J116").Value = 7035.625
$ws.Range("L116").Value = 7035.625
$ws.Range("N116").Value = -13919.625
$ws.Range("H129").Value = 1214.05
$ws.Range("I129").Value = 680.3333
$ws.Range("K129").Value = 2040.9999
$ws.Range("M129").Value = 2959.0001
$ws.Range("H132").Value = 1479.7869
$ws.Range("I132").Value = 1457.2759
$ws.Range("K132").Value = 4371.8277
$ws.Range("M132").Value = -1841.8277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2696.7666
$ws.Range("I61").Value = 1453.6666
$ws.Range("J61").Value = 7669.1665
$ws.Range("K61").Value = 1453.6666
$ws.Range("L61").Value = 7669.1665
$ws.Range("M61").Value = -1241.6666
$ws.Range("N61").Value = -8093.1665
$ws.Range("H132").Value = 4177.109
$ws.Range("I132").Value = 2499
$ws.Range("J132").Value = 8436.923000000001
$ws.Range("K132").Value = 7497
$ws.Range("L132").Value = 25310.769
$ws.Range("M132").Value = -4967
$ws.Range("N132").Value = -30370.769
$ws.Range("H136").Value = 2696.7666
$ws.Range("I136").Value = 1453.6666
$ws.Range("J136").Value = 7669.1665
$ws.Range("K136").Value = 4360.9998
$ws.Range("L136").Value = 23007.4995
$ws.Range("M136").Value = -1810.9998
$ws.Range("N136").Value = -28107.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4279.4287
$ws.Range("I134").Value = 1907.1
$ws.Range("K134").Value = 5721.299999999999
$ws.Range("M134").Value = -3186.299999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4919.61
$ws.Range("I31").Value = 2699.05
$ws.Range("J31").Value = 7034.4287
$ws.Range("K31").Value = 2699.05
$ws.Range("L31").Value = 7034.4287
$ws.Range("M31").Value = -2404.05
$ws.Range("N31").Value = -7624.4287
$ws.Range("H34").Value = 4919.61
$ws.Range("I34").Value = 2699.05
$ws.Range("J34").Value = 7034.4287
$ws.Range("K34").Value = 2699.05
$ws.Range("L34").Value = 7034.4287
$ws.Range("M34").Value = -2497.05
$ws.Range("N34").Value = -7438.4287
$ws.Range("H86").Value = 75782390
$ws.Range("I86").Value = 30338372
$ws.Range("J86").Value = 166670420
$ws.Range("K86").Value = 30338372
$ws.Range("L86").Value = 166670420
$ws.Range("M86").Value = -30337249
$ws.Range("N86").Value = -166672666
$ws.Range("H89").Value = 75782390
$ws.Range("I89").Value = 30338372
$ws.Range("J89").Value = 166670420
$ws.Range("K89").Value = 151691860
$ws.Range("L89").Value = 833352100
$ws.Range("M89").Value = -151686244
$ws.Range("N89").Value = -833363332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2996.75
$ws.Range("I5").Value = 739
$ws.Range("K5").Value = 2217
$ws.Range("M5").Value = -2105
$ws.Range("H107").Value = 20001306
$ws.Range("I107").Value = 690
$ws.Range("J107").Value = 25001460
$ws.Range("K107").Value = 2070
$ws.Range("L107").Value = 75004380
$ws.Range("M107").Value = -150
$ws.Range("N107").Value = -75008220
$ws.Range("H113").Value = 7748.364
$ws.Range("J113").Value = 8423.200000000001
$ws.Range("L113").Value = 25269.6
$ws.Range("N113").Value = -29609.6
$ws.Range("H135").Value = 2996.75
$ws.Range("I135").Value = 739
$ws.Range("K135").Value = 6651
$ws.Range("M135").Value = -4116

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4205321
$ws.Range("I122").Value = 7937706.5
$ws.Range("J122").Value = 6386.875
$ws.Range("K122").Value = 23813119.5
$ws.Range("L122").Value = 19160.625
$ws.Range("M122").Value = -23810669.5
$ws.Range("N122").Value = -24060.625
$ws.Range("H126").Value = 45459212
$ws.Range("I126").Value = 166668670
$ws.Range("J126").Value = 5666.625
$ws.Range("K126").Value = 500006010
$ws.Range("L126").Value = 16999.875
$ws.Range("M126").Value = -500003540
$ws.Range("N126").Value = -21939.875
$ws.Range("H132").Value = 2280.375
$ws.Range("I132").Value = 1245.5385
$ws.Range("J132").Value = 6764.6665
$ws.Range("K132").Value = 3736.6155
$ws.Range("L132").Value = 20293.9995
$ws.Range("M132").Value = -1206.6155
$ws.Range("N132").Value = -25353.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 11203494
$ws.Range("I46").Value = 5749792
$ws.Range("J46").Value = 15878095
$ws.Range("K46").Value = 5749792
$ws.Range("L46").Value = 15878095
$ws.Range("M46").Value = -5749604
$ws.Range("N46").Value = -15878471
$ws.Range("H122").Value = 2734.963
$ws.Range("I122").Value = 2210.587
$ws.Range("J122").Value = 5750.125
$ws.Range("K122").Value = 6631.761
$ws.Range("L122").Value = 17250.375
$ws.Range("M122").Value = -4181.761
$ws.Range("N122").Value = -22150.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 56515
$ws.Range("J119").Value = 56515
$ws.Range("L119").Value = 56515
$ws.Range("N119").Value = -66191
$ws.Range("H122").Value = 15754103
$ws.Range("I122").Value = 21915326
$ws.Range("J122").Value = 8755
$ws.Range("K122").Value = 65745978
$ws.Range("L122").Value = 26265
$ws.Range("M122").Value = -65743528
$ws.Range("N122").Value = -31165
$ws.Range("H126").Value = 52635224
$ws.Range("I126").Value = 100003540
$ws.Range("K126").Value = 300010620
$ws.Range("M126").Value = -300008150

Write-Output "Applied 160 cell updates across 8 sheets"